$wb = $excel.ActiveWorkbook

$wsTopics    = $wb.Worksheets.Item("Topics")
$wsMaterials = $wb.Worksheets.Item("Materials")
$wsWebLogic  = $wb.Worksheets.Item("WebLogic")
$wsGlassFish = $wb.Worksheets.Item("GlassFish")

# ---------------------------------------------------------------------------
# Materials sheet: new "Weblogic" section with two reference links
# ---------------------------------------------------------------------------

# Row 2: section header "Weblogic", styled/merged like the other section
# headers on this workbook (copy format from an existing header, e.g.
# GlassFish!A2:B2, then merge + set the text).
$wsGlassFish.Range("A2:B2").Copy()
$wsMaterials.Range("A2:B2").PasteSpecial(-4122)
$wsMaterials.Range("A2:B2").Merge()
$wsMaterials.Range("A2").Value = "Weblogic"

# Shared-string order in the target workbook has the "Configuring..." /
# CLUST629 pair registered before the "Introduction..." / COHDG5163 pair,
# so write row 4's values before row 3's to reproduce that order, even
# though row 3 precedes row 4 on the sheet.

# Row 4: Configuring and Managing Coherence Clusters
$wsMaterials.Range("A4").Value = "Configuring and Managing Coherence Clusters"
$wsMaterials.Range("B4").Value = "https://docs.oracle.com/middleware/1212/wls/CLUST/coherence.htm#CLUST629"
$wsMaterials.Hyperlinks.Add($wsMaterials.Range("B4"), "https://docs.oracle.com/middleware/1212/wls/CLUST/coherence.htm", "CLUST629")
$wsMaterials.Range("B4").Style = "Hyperlink"

# Row 3: Introduction to Coherence Clusters
$wsMaterials.Range("A3").Value = "Introduction to Coherence Clusters"
$wsMaterials.Range("B3").Value = "https://docs.oracle.com/cd/E18686_01/coh.37/e18677/cluster_overview.htm#COHDG5163"
$wsMaterials.Hyperlinks.Add($wsMaterials.Range("B3"), "https://docs.oracle.com/cd/E18686_01/coh.37/e18677/cluster_overview.htm", "COHDG5163")
$wsMaterials.Range("B3").Style = "Hyperlink"

# ---------------------------------------------------------------------------
# WebLogic sheet: move the frozen-pane selection to the new merged header
# ---------------------------------------------------------------------------
$wsWebLogic.Range("A2:B2").Select()

# ---------------------------------------------------------------------------
# Materials becomes the active sheet/tab, with A3 selected
# ---------------------------------------------------------------------------
$wsMaterials.Activate()
$wsMaterials.Range("A3").Select()
